$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Means" -- add "Within 5 miles" (col F) and "Within 10 miles" (col G)
# of HFC production facility, and refresh the Total Cancer Risk / Total
# Respiratory rows to the recomputed values that came with the new radii.
# ---------------------------------------------------------------------------
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("F1").Value = "Within 5 miles of HFC production facility"
$wsMeans.Range("G1").Value = "Within 10 miles of HFC production facility"

$wsMeans.Range("F2").Value = 97
$wsMeans.Range("G2").Value = 90

$wsMeans.Range("F3").Value = 0.22
$wsMeans.Range("G3").Value = 6.3

$wsMeans.Range("F4").Value = 2.4
$wsMeans.Range("G4").Value = 3.9

$wsMeans.Range("F5").Value = 0.89
$wsMeans.Range("G5").Value = 0.89

$wsMeans.Range("F6").Value = 44
$wsMeans.Range("G6").Value = 48

$wsMeans.Range("F7").Value = 12
$wsMeans.Range("G7").Value = 9.1

$wsMeans.Range("F8").Value = 6.2
$wsMeans.Range("G8").Value = 9.2

# Total Cancer Risk (per million) row -- recomputed means
$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 27
$wsMeans.Range("D9").Value = 30
$wsMeans.Range("E9").Value = 30
$wsMeans.Range("F9").Value = 30
$wsMeans.Range("G9").Value = 31

# Total Respiratory (hazard quotient) row -- recomputed means
$wsMeans.Range("B10").Value = 0.32
$wsMeans.Range("C10").Value = 0.32
$wsMeans.Range("D10").Value = 0.4
$wsMeans.Range("E10").Value = 0.4
$wsMeans.Range("F10").Value = 0.4
$wsMeans.Range("G10").Value = 0.39

# ---------------------------------------------------------------------------
# Sheet "Standard Deviations" -- same two new columns, plus refreshed SD rows
# ---------------------------------------------------------------------------
$wsSD = $wb.Worksheets.Item("Standard Deviations")

$wsSD.Range("F1").Value = "Within 5 mile of HFC production facility SD"
$wsSD.Range("G1").Value = "Within 10 mile of HFC production facility SD"

$wsSD.Range("F2").Value = 2.5
$wsSD.Range("G2").Value = 14

$wsSD.Range("F3").Value = 0.57
$wsSD.Range("G3").Value = 11

$wsSD.Range("F4").Value = 2.4
$wsSD.Range("G4").Value = 4.4

$wsSD.Range("F5").Value = 1.1
$wsSD.Range("G5").Value = 1.7

$wsSD.Range("F6").Value = 9.9
$wsSD.Range("G6").Value = 21

$wsSD.Range("F7").Value = 10
$wsSD.Range("G7").Value = 8.9

$wsSD.Range("F8").Value = 4.5
$wsSD.Range("G8").Value = 9.1

# Total Cancer Risk (per million) SD row -- recomputed
$wsSD.Range("B9").Value = 8.6
$wsSD.Range("C9").Value = 5.7
$wsSD.Range("D9").Value = 0
$wsSD.Range("E9").Value = 0
$wsSD.Range("F9").Value = 0
$wsSD.Range("G9").Value = 2.8

# Total Respiratory (hazard quotient) SD row -- recomputed
$wsSD.Range("B10").Value = 0.14
$wsSD.Range("C10").Value = 0.066
$wsSD.Range("D10").Value = 0
$wsSD.Range("E10").Value = 0.000000000000000028
$wsSD.Range("F10").Value = 0.00000000000000002
$wsSD.Range("G10").Value = 0.027
